# Update "想去人数" (interest count) figures in column F across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets,
# as scraped/regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$ws1.Range("F3").Value = 4887
$ws1.Range("F4").Value = 214
$ws1.Range("F5").Value = 176
$ws1.Range("F6").Value = 131
$ws1.Range("F9").Value = 787
$ws1.Range("F10").Value = 245
$ws1.Range("F11").Value = 1251
$ws1.Range("F12").Value = 133
$ws1.Range("F14").Value = 211
$ws1.Range("F19").Value = 4214
$ws1.Range("F20").Value = 6514
$ws1.Range("F24").Value = 548
$ws1.Range("F26").Value = 4033
$ws1.Range("F28").Value = 55
$ws1.Range("F29").Value = 32
$ws1.Range("F30").Value = 2629
$ws1.Range("F32").Value = 540
$ws1.Range("F34").Value = 317
$ws1.Range("F37").Value = 198
$ws1.Range("F38").Value = 18
$ws1.Range("F40").Value = 990
$ws1.Range("F42").Value = 95
$ws1.Range("F43").Value = 64
$ws1.Range("F44").Value = 509
$ws1.Range("F48").Value = 603

# --- 演出 (sheet2) ---
$ws2.Range("F2").Value = 114

# --- 全部类型 (sheet4) ---
$ws4.Range("F3").Value = 4888
$ws4.Range("F4").Value = 214
$ws4.Range("F5").Value = 176
$ws4.Range("F6").Value = 131
$ws4.Range("F7").Value = 114
$ws4.Range("F10").Value = 787
$ws4.Range("F11").Value = 245
$ws4.Range("F12").Value = 1251
$ws4.Range("F13").Value = 133
$ws4.Range("F15").Value = 211
$ws4.Range("F20").Value = 4214
$ws4.Range("F21").Value = 6514
$ws4.Range("F25").Value = 548
$ws4.Range("F27").Value = 4033
$ws4.Range("F29").Value = 55
$ws4.Range("F30").Value = 32
$ws4.Range("F31").Value = 2629
$ws4.Range("F33").Value = 540
$ws4.Range("F35").Value = 317
$ws4.Range("F38").Value = 198
$ws4.Range("F39").Value = 18
$ws4.Range("F41").Value = 990
$ws4.Range("F43").Value = 95
$ws4.Range("F44").Value = 64
$ws4.Range("F45").Value = 509
$ws4.Range("F49").Value = 603

$wb.Save()
